# Word COM-interop script applying the Mystic Spice Chai Tea market-report edits.
#
# Runtime quirks this script works around:
#  - Once any $d.Tables.Item(N) is accessed, $d.Paragraphs(.Item) stops resolving
#    correctly (it keeps returning the first cell of the touched table for every
#    index). So all paragraph-based edits (the three standalone competitor blurbs)
#    are performed FIRST, before any table is touched.
#  - Storing $d.Tables.Item(1) and $d.Tables.Item(2) in two different variables and
#    then using the older one later resolves to the wrong table (the collection
#    behaves like a single shared cursor). So every table/cell access below re-reads
#    $d.Tables.Item(N) fresh immediately before use instead of caching across calls
#    that touch a different table index.
#  - Find.Execute scoped to a Range taken directly from Cell.Range/Paragraph.Range
#    still searches/replaces against the whole document. Wrapping it as
#    $d.Range(rangeObj.Start, rangeObj.End) makes Find.Execute actually respect the
#    boundaries, so wdReplaceOne only touches the intended occurrence (several of
#    these sentences are substrings of other sentences elsewhere in the document).

$d = $word.ActiveDocument

function Fix-Range-Text($rangeObj, [string]$oldText, [string]$newText) {
    $scoped = $d.Range($rangeObj.Start, $rangeObj.End)
    $ok = $scoped.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 1)
    if (-not $ok) {
        throw "Find/Replace failed for: $oldText"
    }
}

# ---------- Phase 1: standalone paragraphs (competitor descriptions) ----------
# (Must run before any Tables.Item() access -- see note above.)

$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t.Length -ge 8 -and $t.Substring(0, 8) -eq "Teavana:") {
        Fix-Range-Text $p.Range `
            "Teavana: Teavana es una empresa de té basada en Estados Unidos que es propiedad de Starbucks y opera en varios países latinoamericanos, como México, Colombia y Perú." `
            "Teavana: Teavana es una empresa de té de Estados Unidos que es propiedad de Starbucks y opera en varios países latinoamericanos, como México, Colombia y Perú."
        $found = $true
        break
    }
}
if (-not $found) { throw "Teavana paragraph not found" }

$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t.Length -ge 11 -and $t.Substring(0, 11) -eq "Té de David") {
        Fix-Range-Text $p.Range `
            "Té de David: David's Tea es una compañía canadiense de té que tiene presencia en algunos países latinoamericanos, como Chile y Costa Rica." `
            "David's Tea: David's Tea es una empresa de té canadiense que tiene presencia en algunos países latinoamericanos, como Chile y Costa Rica."
        $found = $true
        break
    }
}
if (-not $found) { throw "David's Tea paragraph not found" }

$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t.Length -ge 15 -and $t.Substring(0, 15) -eq "Marcas locales:") {
        Fix-Range-Text $p.Range `
            "Marcas locales: También hay varias marcas locales que ofrecen productos de té Chai en América Latina, como Mate Factor, Chai Mate y Chai Brasil." `
            "Marcas locales: También hay varias marcas locales que ofrecen productos de té chai en América Latina, como Mate Factor, Chai Mate y Chai Brasil."
        $found = $true
        break
    }
}
if (-not $found) { throw "Marcas locales paragraph not found" }

# ---------- Phase 2: Table 1 (product description table) ----------

# Header cell (1,2): "Descripcion del producto" should be bold like its neighbour.
$cell = $d.Tables.Item(1).Cell(1, 2)
$boldRng = $d.Range($cell.Range.Start, $cell.Range.End)
$boldRng.Font.Bold = 1

# (2,1): "te chai premium Mystic Spice" -> "Te chai premium Mystic Spice"
$cell = $d.Tables.Item(1).Cell(2, 1)
Fix-Range-Text $cell.Range `
    "té chai premium Mystic Spice" `
    "Té chai premium Mystic Spice"

# (4,1): authentic blend sentence
$cell = $d.Tables.Item(1).Cell(4, 1)
Fix-Range-Text $cell.Range `
    "Mezcla auténtica: Nuestra chai es una mezcla armónica de hojas de té negro premium y una selección de especias molidas, incluyendo canela, cardamomo, cloves, jengibre y pimienta negra." `
    "Mezcla auténtica: nuestro chai es una mezcla armoniosa de hojas de té negro premium y una selección de especias molidas, incluyendo canela, cardamomo, clavo, jengibre y pimienta negra."

# (4,2): two separate runs both change
$cell = $d.Tables.Item(1).Cell(4, 2)
Fix-Range-Text $cell.Range `
    "Ingredientes de mejora de la salud: Cada ingrediente del Místico Spice Chai Tea se elige para sus beneficios naturales para la salud." `
    "Ingredientes beneficiosos para la salud: cada ingrediente del té chai Mystic Spice se elige para sus beneficios naturales para la salud."

$cell = $d.Tables.Item(1).Cell(4, 2)
Fix-Range-Text $cell.Range `
    "El jengibre y el cardamomo ayudan a la digestión, la canela ayuda a regular el azúcar en sangre y los clavos agregan un impulso de antioxidantes." `
    "El jengibre y el cardamomo ayudan a la digestión, la canela ayuda a regular el azúcar en sangre y el clavo aumento los antioxidantes."

# (5,1): aroma/flavor sentence
$cell = $d.Tables.Item(1).Cell(5, 1)
Fix-Range-Text $cell.Range `
    "Rico Aroma y Sabor: El aroma cálido, picante y profundo, vigorizante sabor de nuestra chai hacen que sea la bebida perfecta para comenzar su día o relajarse por la noche." `
    "Aroma y sabor intensos: el aroma cálido y especiado, y el sabor profundo y vigorizante de nuestro chai hacen que sea la bebida perfecta para comenzar el día o relajarse por la noche."

# (5,2): preparation options sentence
$cell = $d.Tables.Item(1).Cell(5, 2)
Fix-Range-Text $cell.Range `
    "Opciones versátiles de preparación: Ya sea que amas tu chai vaporing caliente, como un refrescante té helado, o como una latte cremosa, nuestra mezcla es lo suficientemente versátil como para adaptarte a cualquier preferencia." `
    "Opciones versátiles de preparación: Ya sea que le guste su chai bien caliente, como un refrescante té helado, o como una latte cremoso, nuestra mezcla es lo suficientemente versátil como para adaptarse a cualquier preferencia."

# (6,1): sustainable origin sentence
$cell = $d.Tables.Item(1).Cell(6, 1)
Fix-Range-Text $cell.Range `
    "Origen sostenible: Comprometidos con la sostenibilidad, originamos nuestros ingredientes de granjas a pequeña escala que practican la agricultura ecológica, garantizando no sólo la mejor calidad, sino también el bienestar de nuestro planeta." `
    "Origen sostenible: al estar comprometidos con la sostenibilidad, obtenemos nuestros ingredientes de pequeñas explotaciones que practican la agricultura ecológica, garantizando no solo la mejor calidad, sino también el bienestar de nuestro planeta."

# (6,2): elegant packaging sentence
$cell = $d.Tables.Item(1).Cell(6, 2)
Fix-Range-Text $cell.Range `
    "Empaquetado elegante: El té de Spice Chai místico viene en un empaquetado elegante, ecológico, lo que lo convierte en un regalo ideal para los amantes del té o un lujoso trato para usted mismo." `
    "Envase elegante: el té chai Mystic Spice viene en un envase elegante, ecológico, lo que lo convierte en el regalo ideal para los amantes del té o un capricho lujoso para ti mismo."

# (7,1): customer satisfaction guarantee sentence
$cell = $d.Tables.Item(1).Cell(7, 1)
Fix-Range-Text $cell.Range `
    "Garantía de satisfacción del cliente: Estamos detrás de nuestro producto y ofrecemos una garantía de satisfacción." `
    "Garantía de satisfacción del cliente: Respaldamos nuestro producto y ofrecemos una garantía de satisfacción."

# (7,2): "Ideal para" sentence
$cell = $d.Tables.Item(1).Cell(7, 2)
Fix-Range-Text $cell.Range `
    "Ideal para: entusiastas del té, individuos conscientes de la salud, amantes de bebidas calientes, especiadas, y cualquier persona que busca explorar los ricos sabores de la chai india tradicional." `
    "Ideal para: los apasionados del té, las personas conscientes de la salud, los amantes de las bebidas calientes especiadas, y cualquier persona que desea explorar los sabores intensos del chai indio tradicional."

# ---------- Phase 3: Table 2 (competitor table) ----------

# Header row bold fixes: Marca / Cuota de mercado (%) / Estrategia de precios
$cell = $d.Tables.Item(2).Cell(1, 1)
$boldRng = $d.Range($cell.Range.Start, $cell.Range.End)
$boldRng.Font.Bold = 1

$cell = $d.Tables.Item(2).Cell(1, 2)
$boldRng = $d.Range($cell.Range.Start, $cell.Range.End)
$boldRng.Font.Bold = 1

$cell = $d.Tables.Item(2).Cell(1, 3)
$boldRng = $d.Range($cell.Range.Start, $cell.Range.End)
$boldRng.Font.Bold = 1

# (1,3) text: "Estrategia de precios" -> "Estrategias de precios"
$cell = $d.Tables.Item(2).Cell(1, 3)
Fix-Range-Text $cell.Range `
    "Estrategia de precios" `
    "Estrategias de precios"

Write-Host "All edits applied."
